$p = $ppt.ActivePresentation

# --- Slide 10: "Megvalósítás" progress table -------------------------------
# Update several cell values (team member workload table) and switch the
# table style to the new style GUID.
$s10 = $p.Slides.Item(10)
$tbl10 = $s10.Shapes.Item(3).Table

$tbl10.Cell(4, 6).Shape.TextFrame.TextRange.Text = "42"

$tbl10.Cell(9, 2).Shape.TextFrame.TextRange.Text = "28"
$tbl10.Cell(9, 3).Shape.TextFrame.TextRange.Text = "20"
$tbl10.Cell(9, 4).Shape.TextFrame.TextRange.Text = "20"
$tbl10.Cell(9, 6).Shape.TextFrame.TextRange.Text = "83"

$tbl10.Cell(10, 2).Shape.TextFrame.TextRange.Text = "99"
$tbl10.Cell(10, 3).Shape.TextFrame.TextRange.Text = "80"
$tbl10.Cell(10, 4).Shape.TextFrame.TextRange.Text = "80"
$tbl10.Cell(10, 6).Shape.TextFrame.TextRange.Text = "343"

$tbl10.ApplyStyle("{30BCC980-B1AF-4FEF-85EA-2FC768ABBB76}")

# --- Slide 12: two code-snippet tables --------------------------------------
$s12 = $p.Slides.Item(12)
$s12.Shapes.Item(4).Table.ApplyStyle("{30BCC980-B1AF-4FEF-85EA-2FC768ABBB76}")
$s12.Shapes.Item(5).Table.ApplyStyle("{30BCC980-B1AF-4FEF-85EA-2FC768ABBB76}")

# --- Slide 13: code-snippet table -------------------------------------------
$s13 = $p.Slides.Item(13)
$s13.Shapes.Item(4).Table.ApplyStyle("{30BCC980-B1AF-4FEF-85EA-2FC768ABBB76}")

# --- Slide 17: large results table ------------------------------------------
$s17 = $p.Slides.Item(17)
$s17.Shapes.Item(3).Table.ApplyStyle("{38192677-2D9B-4A0D-92FA-7B8146FB3DC1}")
